$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text updates -------------------------------------------------
# G1: "Days (e.g., Monday, Wednesday)" -> "Day (e.g., Monday)"
$ws.Range("G1").Value = "Day (e.g., Monday)"

# H1: " Room Type" -> " Room Type(Lecture, Laboratory)"
$ws.Range("H1").Value = " Room Type(Lecture, Laboratory)"

# --- Remove the sample data rows (rows 2-4) ---------------------------------
$ws.Range("A2:H4").ClearContents()

# --- Add a list data-validation on column H (room type) ---------------------
# xlValidateList=3, xlValidAlertStop=1 -> allowBlank/showInput/showError all on
$rng = $ws.Range("H2:H1048576")
$rng.Validation.Add(3, 1, 1, '"Lecture, Laboratory"')

# --- Sheet view: reset horizontal scroll + move selection to A2 -------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
[void]$ws.Range("A2").Select()
